$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '53.457.20'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.65%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.157.75'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.64%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '397.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.25%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.43'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.43%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.544'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.05%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.02%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.609'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.12%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.46'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.26%  '

# Row 11
$ws.Range("E11").Value = '  +1.05%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0868'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.68%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.645.42'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.33%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.89'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.27%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.93'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.14%  '

# Row 16
$ws.Range("E16").Value = '  +7.79%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.159.30'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.66%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.69'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.49%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '53.190.67'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.04%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.76%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.88'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.59%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0977'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.31%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.35'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.66%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '271.02'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.97%  '

# Row 25
$ws.Range("E25").Value = '  +1.63%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.10'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.37%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '27.65'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.88%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.52'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.23%  '

# Row 29
$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.171'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.34%  '

# Row 30
$ws.Range("B30").Value = 'Dai'
$ws.Range("C30").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.04%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.109'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.50%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.93'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.31%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '36.93'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.35%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0492'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +10.71%  '

# Row 35
$ws.Range("E35").Value = '  +0.73%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '50.54'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.10%  '

# Row 37
$ws.Range("E37").Value = '  -0.06%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.49'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.04%  '

# Row 39
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.83'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +11.03%  '

# Row 40
$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.06'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.58%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.292'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.29%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '17.22'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.32%  '

# Row 43
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.90'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.22%  '

# Row 44
$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '130.65'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.96%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.117'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.63%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.32'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.99%  '

# Row 47
$ws.Range("E47").Value = '  -1.45%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.07'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.60%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.086.95'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.46%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0526'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +32.88%  '

# Row 51
$ws.Range("B51").Value = 'BEAM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0331'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.85%  '
